$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# The worksheet is protected (legacy password hash) and every cell's style
# is implicitly "locked" (Excel default). Calling Range.Unprotect()/Protect()
# on the sheet rewrites the <sheetProtection> element to a brand-new
# SHA-512 hash, which would wrongly show up as a change. Instead, toggle the
# .Locked flag of only the specific cells we need to touch: unlocking and
# re-locking a cell is allowed even while the sheet stays protected, and it
# leaves <sheetProtection> completely untouched.

function Set-CellValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.Locked = $false
    $rng.Value2 = $value
    $rng.Locked = $true
}

# --- Update the "as of" date in the confidentiality / disclaimer text (A7) ---
$a7 = $ws.Range("A7")
$a7.Locked = $false
$oldText = $a7.Value2
$newText = $oldText -replace "2021-05-17", "2021-05-18"
$a7.Value2 = $newText
$a7.Locked = $true

# --- Update the Weight / Percent Change figures for rows 2-4 ---
Set-CellValue "D2" 0.8491534485831403
Set-CellValue "E2" 0.002140788313814479

Set-CellValue "D3" 0.1508465514168598
Set-CellValue "E3" 0.01342155009451806

Set-CellValue "E4" 0.003842452325788948
